# Sync attendance_reports: reorder the "Recorded By" (column G) author
# lists on the "Session Analysis Results" sheet so "System" sorts after
# the other recorder name(s) instead of before them, e.g.:
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System, system" -> "backup@backdoor.com, system, System"
# All other "Recorded By" values (single recorder, or lists that already
# have "System" last) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
        $changed++
    }
    elseif ($text -eq "backup@backdoor.com, System, system") {
        $cell.Value = "backup@backdoor.com, system, System"
        $changed++
    }
}

Write-Output "Reordered Recorded By values in $changed cell(s)."
